$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Locate the paragraph that begins "Please read Using Pressure
#    Canners before beginning...." and the lone "\n" paragraph that
#    immediately follows it (the blank-line separator before the
#    "Procedure:" paragraph).
# ------------------------------------------------------------------
$introIndex = -1
$idx = 0
foreach ($p in $d.Paragraphs) {
    $idx = $idx + 1
    $t = $p.Range.Text
    if ($t.StartsWith("Please read Using Pressure Canners before beginning")) {
        $introIndex = $idx
    }
}

$introPara  = $d.Paragraphs.Item($introIndex)
$blankPara  = $d.Paragraphs.Item($introIndex + 1)

# Delete the blank "\n" paragraph that sits between the intro
# paragraph and the "Procedure: ..." paragraph, merging the
# following paragraph's content upward.
$blankPara.Range.Delete()

# ------------------------------------------------------------------
# 2. Replace the intro paragraph's text with "Procedure" and make the
#    trailing literal "\n" its own run (two runs total), matching the
#    target paragraph structure.
# ------------------------------------------------------------------
$introPara = $d.Paragraphs.Item($introIndex)
$bodyRange = $d.Range($introPara.Range.Start, $introPara.Range.End - 1)
$bodyRange.Text = "Procedure"

$introPara = $d.Paragraphs.Item($introIndex)
$tailInsertionPoint = $d.Range($introPara.Range.End - 1, $introPara.Range.End - 1)
$tailInsertionPoint.InsertAfter("\n")

# Force the newly inserted "\n" text to live in its own run (instead
# of being silently merged back into the "Procedure" run) by toggling
# a character property on just that span.
$introPara = $d.Paragraphs.Item($introIndex)
$tailRunRange = $d.Range($introPara.Range.End - 3, $introPara.Range.End - 1)
$tailRunRange.Bold = 1
$tailRunRange.Bold = 0

# ------------------------------------------------------------------
# 3. Strip the leading "Procedure: " label from the paragraph that
#    used to read "Procedure: Wash beans and trim ends. ..." so it
#    now simply starts with "Wash beans...". The paragraph's second
#    run (the lone "n" that together with the trailing "\" forms the
#    literal "\n") must remain a separate run.
# ------------------------------------------------------------------
$procIndex = -1
$idx = 0
foreach ($p in $d.Paragraphs) {
    $idx = $idx + 1
    $t = $p.Range.Text
    if ($t.StartsWith("Procedure: Wash beans and trim ends")) {
        $procIndex = $idx
    }
}

$procPara = $d.Paragraphs.Item($procIndex)
$prefix = "Procedure: "
$prefixRange = $d.Range($procPara.Range.Start, $procPara.Range.Start + $prefix.Length)
$prefixRange.Delete()

# Re-establish the run boundary between "...pieces.\" and the final
# "n" (otherwise the save step silently merges the two identically
# formatted runs back together).
$procPara = $d.Paragraphs.Item($procIndex)
$finalCharRange = $d.Range($procPara.Range.End - 2, $procPara.Range.End - 1)
$finalCharRange.Bold = 1
$finalCharRange.Bold = 0

Write-Host "Done."
